$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Smeta")

foreach ($source in $wb.LinkSources(1)) {
    $wb.BreakLink($source, 1)
}

$values = @(
    "Clay Paky Alpha Beam 1500",
    "Clay Paky Alpha Profile 1500",
    "Skylight F230",
    "Esdelumen Smart p6.0 0.576 x0.576m, 96x96 pixels, 0.25 sq.m",
    "Ферма  треугольная/truss triangle 30х30, 2м (стрелы)",
    "кубы 30х30"
)

$row = 10
foreach ($v in $values) {
    $ws.Cells.Item($row, 3).Value = $v
    $row = $row + 1
}

$ws.Range("D18").Select()

$wb.Save()
